# Update Jupyter Notebook Examples
# Append newer "Kiel Arbeitslose" yearly figures (2017-2020) to the data
# table on Sheet1, matching the formatting used for the other data rows,
# and update the sheet's view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New yearly rows to append after the existing data (last row = 28, year 2016)
$newRows = @(
    @(2017, 11939),
    @(2018, 10260),
    @(2019, 10009),
    @(2020, 11776)
)

$startRow = 29
$r = $startRow
foreach ($pair in $newRows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    # Match styling applied to the new "Arbeitslose" values (black font color)
    $ws.Cells.Item($r, 2).Font.Color = 0
    $r++
}

# Update the view: scroll so row 12 is at the top, and move the active
# selection to a single cell further down in the sheet.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
